$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Clear rows 1 and 2 entirely (B1:E2), which held summary stats that are no
# longer needed.
$ws.Range("B1:E2").Clear()

# Remove the "nothing" column (E) from the header row; data now only spans
# B:D.
$ws.Range("E3").Clear()

# Update the active selection to match the new focus cell.
$ws.Range("E8").Select()
